$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1734.5
$ws.Range("J17").Value = 2040.8
$ws.Range("L17").Value = 6122.4
$ws.Range("N17").Value = -6458.4

# Sheet ALC (index 1), row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 748.75
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 665
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 665
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -1317

# Sheet ALC (index 1), row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 318.13635
$ws.Range("I33").Value = 238.88889
$ws.Range("K33").Value = 238.88889
$ws.Range("M33").Value = -9.888890000000004

# Sheet ALC (index 1), row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 12000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 12000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 12000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -12812

# Sheet ALC (index 1), row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 12000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 12000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 12000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -14808

# Sheet ALC (index 1), row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1269.1904
$ws.Range("I92").Value = 493.625
$ws.Range("J92").Value = 3751
$ws.Range("K92").Value = 493.625
$ws.Range("L92").Value = 3751
$ws.Range("M92").Value = 754.375
$ws.Range("N92").Value = -6247

# Sheet ALC (index 1), row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 655.9167
$ws.Range("I125").Value = 409.5
$ws.Range("K125").Value = 3685.5
$ws.Range("M125").Value = -1225.5

# Sheet ALC (index 1), row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 44974.258
$ws.Range("I132").Value = 45934.83
$ws.Range("K132").Value = 137804.49
$ws.Range("M132").Value = -135274.49

# Sheet ALC (index 1), row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3398.5676
$ws.Range("J137").Value = 3067.92
$ws.Range("L137").Value = 9203.76
$ws.Range("N137").Value = -14303.76

# Sheet ALC (index 1), row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4624.514
$ws.Range("I138").Value = 7500
$ws.Range("K138").Value = 22500
$ws.Range("M138").Value = -17360

# Sheet ARM (index 2), row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1095.909
$ws.Range("I110").Value = 956.875
$ws.Range("K110").Value = 956.875
$ws.Range("M110").Value = 1088.125

# Sheet BSM (index 3), row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 11322
$ws.Range("I99").Value = 5511.737
$ws.Range("K99").Value = 5511.737
$ws.Range("M99").Value = -4013.737

# Sheet CRP (index 4), row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 111499.6
$ws.Range("I31").Value = 171005.11
$ws.Range("K31").Value = 171005.11
$ws.Range("M31").Value = -170710.11

# Sheet CRP (index 4), row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 111499.6
$ws.Range("I34").Value = 171005.11
$ws.Range("K34").Value = 171005.11
$ws.Range("M34").Value = -170803.11

# Sheet CRP (index 4), row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 14281
$ws.Range("J88").Value = 14281
$ws.Range("L88").Value = 14281
$ws.Range("N88").Value = -15093

# Sheet CRP (index 4), row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 14281
$ws.Range("J91").Value = 14281
$ws.Range("L91").Value = 14281
$ws.Range("N91").Value = -17089

# Sheet CUL (index 5), row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6038.5293
$ws.Range("I75").Value = 2542.7144
$ws.Range("J75").Value = 8485.6
$ws.Range("K75").Value = 7628.1432
$ws.Range("L75").Value = 25456.8
$ws.Range("M75").Value = -6630.1432
$ws.Range("N75").Value = -27452.8

# Sheet CUL (index 5), row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 6038.5293
$ws.Range("I78").Value = 2542.7144
$ws.Range("J78").Value = 8485.6
$ws.Range("K78").Value = 22884.4296
$ws.Range("L78").Value = 76370.40000000001
$ws.Range("M78").Value = -17892.4296
$ws.Range("N78").Value = -86354.40000000001

# Sheet GSM (index 6), row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 275.7619
$ws.Range("I2").Value = 179.07692
$ws.Range("J2").Value = 432.875
$ws.Range("K2").Value = 179.07692
$ws.Range("L2").Value = 432.875
$ws.Range("M2").Value = -66.07692
$ws.Range("N2").Value = -658.875

# Sheet GSM (index 6), row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1538.2307
$ws.Range("I97").Value = 1833
$ws.Range("K97").Value = 1833
$ws.Range("M97").Value = -1337

# Sheet GSM (index 6), row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 83123.71000000001
$ws.Range("I122").Value = 129318.625
$ws.Range("J122").Value = 21530.5
$ws.Range("K122").Value = 387955.875
$ws.Range("L122").Value = 64591.5
$ws.Range("M122").Value = -385505.875
$ws.Range("N122").Value = -69491.5

# Sheet LTW (index 7), row 3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2004
$ws.Range("I3").Value = 2004
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2004
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1892
$ws.Range("N3").ClearContents()

# Sheet LTW (index 7), row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2830.35
$ws.Range("I7").Value = 1966.8334
$ws.Range("J7").Value = 10602
$ws.Range("K7").Value = 1966.8334
$ws.Range("L7").Value = 10602
$ws.Range("M7").Value = -1854.8334
$ws.Range("N7").Value = -10826

# Sheet LTW (index 7), row 15
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 2004
$ws.Range("I15").Value = 2004
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2004
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1834
$ws.Range("N15").ClearContents()

# Sheet LTW (index 7), row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1674.75
$ws.Range("I22").Value = 1679.8
$ws.Range("K22").Value = 1679.8
$ws.Range("M22").Value = -1384.8

# Sheet LTW (index 7), row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1674.75
$ws.Range("I27").Value = 1679.8
$ws.Range("K27").Value = 1679.8
$ws.Range("M27").Value = -1572.8

# Sheet LTW (index 7), row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3277.7896
$ws.Range("I40").Value = 2836.4546
$ws.Range("J40").Value = 3884.625
$ws.Range("K40").Value = 2836.4546
$ws.Range("L40").Value = 3884.625
$ws.Range("M40").Value = -2700.4546
$ws.Range("N40").Value = -4156.625

# Sheet LTW (index 7), row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 80782.69500000001
$ws.Range("I122").Value = 3438.8572
$ws.Range("K122").Value = 10316.5716
$ws.Range("M122").Value = -7866.571599999999

# Sheet LTW (index 7), row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2830.35
$ws.Range("I126").Value = 1966.8334
$ws.Range("J126").Value = 10602
$ws.Range("K126").Value = 5900.5002
$ws.Range("L126").Value = 31806
$ws.Range("M126").Value = -3430.5002
$ws.Range("N126").Value = -36746

# Sheet WVR (index 8), row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 260228.28
$ws.Range("I62").Value = 362399.6
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 362399.6
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -361775.6
$ws.Range("N62").Value = -6048

# Sheet WVR (index 8), row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 260228.28
$ws.Range("I65").Value = 362399.6
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 1811998
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -1808878
$ws.Range("N65").Value = -30240

# Sheet WVR (index 8), row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3000.6667
$ws.Range("I122").Value = 2943.7144
$ws.Range("K122").Value = 8831.143199999999
$ws.Range("M122").Value = -6381.143199999999

# Sheet WVR (index 8), row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2166.5293
$ws.Range("I132").Value = 2014.9286
$ws.Range("K132").Value = 6044.7858
$ws.Range("M132").Value = -3514.7858

# Sheet WVR (index 8), row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 99999.5
$ws.Range("J133").Value = 99999.5
$ws.Range("L133").Value = 99999.5
$ws.Range("N133").Value = -110119.5
